# Daily Satellite Data Update
# Shifts the pass-prediction data forward by one day: refines the
# predicted times/percentages for the existing dates (rows 2-5) and
# appends a brand-new row (row 6) for 2026-02-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 2 (20260221--01) - refine predicted times / percentages
# ---------------------------------------------------------------
$ws.Range("C2").Value = "04:01"
$ws.Range("D2").Value = "00:00"
$ws.Range("E2").Value = "05:07:12"
$ws.Range("F2").Value = "05:10:04"
$ws.Range("G2").Value = "05:12:05"
$ws.Range("H2").Value = "05:14:05"
$ws.Range("I2").Value = "05:16:58"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "-"
$ws.Range("M2").Value = "A+B"
$ws.Range("N2").Value = 2
$ws.Range("O2").Value = 99
$ws.Range("P2").Value = 90
$ws.Range("Q2").Value = 95
$ws.Range("R2").Value = 50

# ---------------------------------------------------------------
# Row 3 (20260222--01)
# ---------------------------------------------------------------
$ws.Range("C3").Value = "03:08"
$ws.Range("D3").Value = "03:08"
$ws.Range("E3").Value = "04:28:54"
$ws.Range("F3").Value = "04:32:02"
$ws.Range("G3").Value = "04:33:36"
$ws.Range("H3").Value = "04:35:10"
$ws.Range("I3").Value = "04:38:20"
$ws.Range("J3").Value = "4°"
$ws.Range("K3").Value = "04:30:11"
$ws.Range("M3").Value = "A"
$ws.Range("N3").Value = 4
$ws.Range("P3").Value = 86
$ws.Range("Q3").Value = 76
$ws.Range("R3").Value = 80

# ---------------------------------------------------------------
# Row 4 (20260223--01)
# ---------------------------------------------------------------
$ws.Range("C4").Value = "01:33"
$ws.Range("D4").Value = "01:33"
$ws.Range("E4").Value = "03:50:35"
$ws.Range("F4").Value = "03:54:19"
$ws.Range("G4").Value = "03:55:05"
$ws.Range("H4").Value = "03:55:52"
$ws.Range("I4").Value = "03:59:37"
$ws.Range("J4").Value = "9°"
$ws.Range("K4").Value = "03:53:36"
$ws.Range("L4").Value = -20.7
$ws.Range("M4").Value = "A"
$ws.Range("N4").Value = 4
$ws.Range("P4").Value = 97
$ws.Range("Q4").Value = 92

# ---------------------------------------------------------------
# Row 5 (20260224--01)
# ---------------------------------------------------------------
$ws.Range("C5").Value = "05:51"
$ws.Range("D5").Value = "05:51"
$ws.Range("E5").Value = "04:48:23"
$ws.Range("F5").Value = "04:50:46"
$ws.Range("G5").Value = "04:53:41"
$ws.Range("H5").Value = "04:56:37"
$ws.Range("I5").Value = "04:59:01"
$ws.Range("J5").Value = "8°"
$ws.Range("K5").Value = "04:50:23"
$ws.Range("M5").Value = "A"
$ws.Range("N5").Value = 2
$ws.Range("O5").Value = 100
$ws.Range("P5").Value = 98
$ws.Range("Q5").Value = 45
$ws.Range("R5").Value = 96

# ---------------------------------------------------------------
# Row 6 (20260225--01) - brand new row appended for the new day
# ---------------------------------------------------------------
$ws.Range("A6").Value = "20260225--01"
$ws.Range("B6").Value = 22
$ws.Range("C6").Value = "05:37"
$ws.Range("D6").Value = "04:20"
$ws.Range("E6").Value = "04:09:47"
$ws.Range("F6").Value = "04:12:13"
$ws.Range("G6").Value = "04:15:01"
$ws.Range("H6").Value = "04:17:50"
$ws.Range("I6").Value = "04:20:17"
$ws.Range("J6").Value = "17°"
$ws.Range("K6").Value = "04:13:30"
$ws.Range("L6").Value = -17
$ws.Range("M6").Value = "A+B"
$ws.Range("N6").Value = 2
$ws.Range("O6").Value = 75
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 69

# ---------------------------------------------------------------
# Recolor the O:R "heat-map" cells to reflect the refreshed values
# (colour-scale style fills, blue = low, white/near-white = mid,
# red = flagged column O)
# ---------------------------------------------------------------
$ws.Range("O2").Interior.Color = 8351984
$ws.Range("P2").Interior.Color = 16380910
$ws.Range("Q2").Interior.Color = 16579062
$ws.Range("R2").Interior.Color = 14927274

$ws.Range("O3").Interior.Color = 8351984
$ws.Range("P3").Interior.Color = 16248550
$ws.Range("Q3").Interior.Color = 15852244
$ws.Range("R3").Interior.Color = 16050397

$ws.Range("O4").Interior.Color = 8351984
$ws.Range("P4").Interior.Color = 16579062
$ws.Range("Q4").Interior.Color = 16380910
$ws.Range("R4").Interior.Color = 16777215

$ws.Range("O5").Interior.Color = 8351984
$ws.Range("P5").Interior.Color = 16777215
$ws.Range("Q5").Interior.Color = 14794914
$ws.Range("R5").Interior.Color = 16579062

$ws.Range("O6").Interior.Color = 8351984
$ws.Range("P6").Interior.Color = 13143125
$ws.Range("Q6").Interior.Color = 13143125
$ws.Range("R6").Interior.Color = 15654092

# ---------------------------------------------------------------
# Extend every conditional-formatting rule down to include row 6
# ---------------------------------------------------------------
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($col in $cols) {
    $oldRange = $ws.Range($col + "2:" + $col + "5")
    $newRange = $ws.Range($col + "2:" + $col + "6")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($newRange)
    }
}
